# Edit format tipe rumah
# Adds a "terbilang" (amount spelled out in words) column D to the
# tipe_rumah sheet, formats it like the existing harga_sewa/tipe columns,
# widens the new column, removes the old frozen-pane view, and leaves the
# existing A/B/C data untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cell D1 -- copy the look of the other header cells (C1)
#    then set its text.
# ---------------------------------------------------------------------
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "terbilang"

# ---------------------------------------------------------------------
# 2. New data cells D2:D31 -- copy the look of the harga_sewa column
#    (C2), which carries the border/fill/font formatting, then switch
#    the number format to Text (the column holds words, not numbers)
#    while keeping the right/center alignment.
# ---------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("D2:D31").PasteSpecial(-4122)
$ws.Range("D2:D31").NumberFormat = "@"

$terbilang = @(
    "Satu Juta Lima Ratus Tiga Puluh Satu Ribu Sembilan Ratus Enam Puluh Lima Rupiah",
    "Satu Juta Tiga Ratus Empat Puluh Delapan Ribu Dua Ratus Rupiah",
    "Satu Juta Seratus Enam Puluh Empat Ribu Tiga Ratus Rupiah",
    "Satu Juta Seratus Tiga Ribu Seratus Rupiah",
    "Satu Juta Dua Puluh Tujuh Ribu Lima Ratus Rupiah",
    "Satu Juta Sebelas Ribu Tiga Ratus Rupiah",
    "Sembilan Ratus Sembilan Belas Ribu Dua Ratus Rupiah",
    "Sembilan Ratus Enam Ribu Sembilan Ratus Rupiah",
    "Delapan Ratus Tujuh Puluh Ribu Seratus Lima Puluh Rupiah",
    "Delapan Ratus Lima Puluh Delapan Ribu Rupiah",
    "Delapan Ratus Lima Puluh Delapan Ribu Rupiah",
    "Tujuh Ratus Lima Puluh TIga Ribu Tujuh Ratus Lima Puluh Rupiah",
    "Tujuh Ratus Tiga Puluh Lima Ribu Tiga Ratus Empat Puluh Lima Rupiah",
    "Tujuh Ratus Empat Ribu Tujuh Ratus Rupiah",
    "Enam Ratus Dua Belas Ribu Rupiah",
    "Lima Rates Delapan Puluh Delapan Ribu Rupiah",
    "Empat Ratus Lima Puluh Sembilan Ribu Enam Ratus Rupiah",
    "Empat Ratus Tiga Puluh Sembilan Ribu Dua Ratus Rupiah",
    "Empar Ratus Delapan Ribu Empat Ratus Lima Puluh Rupiah",
    "Tiga Ratus Delapan Puluh Dua Sembilan Ratus Lima Puluh",
    "Tiga Ratus Lima Puluh Tujuh Ribu Empat Ratus Lima Rupiah",
    "Tiga Ratus Tiga Puluh Enam Ribu Sembilan Ratus Rupiah",
    "Tiga Ratus Tiga Puluh Satu Ribu Delapan Ratus Tujuh Puluh Lima Rupiah",
    "Tiga Ratus Enam Ribu Empat Ratus Lima Puluh Rupiah",
    "Dua Ratus Tujuh Puluh Lima Ribu Tujuh Ratus Rupiah",
    "Dua Ratus Lima Puluh Lima Ribu Tiga Ratus Rupiah",
    "Dua Ratus Dua Puluh Sembilan Ribu Delapan Ratus Rupiah",
    "Dua Ratus Empat Ribu Rupiah",
    "Seratus Delapan Puluh Tiga Ribu Sempbilan Ratus Rupiah",
    "Tiga Puluh Enam Ribu Rupiah"
)

$row = 2
foreach ($t in $terbilang) {
    $ws.Cells.Item($row, 4).Value = $t
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. Widen the new terbilang column so the long sentences are readable.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 87.2701142857143

# ---------------------------------------------------------------------
# 4. Drop the frozen header/column pane that used to split the view.
# ---------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
